$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 452 (shifts existing rows 452-462 down to 453-463),
# inheriting the formatting (date number format on column D) from the row above.
$ws.Rows(452).Insert()

# Populate the newly inserted row 452 with the new weekly price entry.
$ws.Cells.Item(452, 1).Value  = 9
$ws.Cells.Item(452, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(452, 3).Value  = "Metropolitana"
$ws.Cells.Item(452, 4).Value  = 45239
$ws.Cells.Item(452, 5).Value  = 13
$ws.Cells.Item(452, 6).Value  = 100112043
$ws.Cells.Item(452, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(452, 8).Value  = "Sin especificar"
$ws.Cells.Item(452, 9).Value  = "Primera"
$ws.Cells.Item(452, 10).Value = 97
$ws.Cells.Item(452, 11).Value = 14000
$ws.Cells.Item(452, 12).Value = 15000
$ws.Cells.Item(452, 13).Value = 14485
$ws.Cells.Item(452, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(452, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(452, 16).Value = 241
$ws.Cells.Item(452, 17).Value = 60
$ws.Cells.Item(452, 18).Value = "Hortaliza"
